# Update the crypto price list (columns B-E) to reflect the latest scrape.
# Rows 2-51 on the active sheet; only the cells that actually changed are
# touched (the diff is value-only: coin name/link swap on rows 31-32, plus
# refreshed Price/Volume(1h) figures everywhere else).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="43.111.92"; E="  -0.53%  " },
    @{ Row=3; D="2.380.44"; E="  +1.02%  " },
    @{ Row=4; D="1.00"; E="  +0.01%  " },
    @{ Row=5; D="303.88"; E="  +0.12%  " },
    @{ Row=6; D="96.50"; E="  +0.63%  " },
    @{ Row=7; D="0.504"; E="  +0.01%  " },
    @{ Row=8; E="  -0.03%  " },
    @{ Row=9; E="  -3.08%  " },
    @{ Row=10; D="34.78"; E="  +1.33%  " },
    @{ Row=11; E="  +3.10%  " },
    @{ Row=12; D="0.0791"; E="  +0.34%  " },
    @{ Row=13; D="18.46"; E="  -0.61%  " },
    @{ Row=14; D="6.82"; E="  +0.55%  " },
    @{ Row=15; D="2.744.38"; E="  +0.84%  " },
    @{ Row=16; D="2.387.21"; E="  +2.21%  " },
    @{ Row=17; D="0.809"; E="  +1.36%  " },
    @{ Row=18; D="43.117.92"; E="  -0.40%  " },
    @{ Row=19; D="6.33"; E="  +1.58%  " },
    @{ Row=20; D="11.99"; E="  -2.09%  " },
    @{ Row=21; E="  -0.13%  " },
    @{ Row=22; D="68.15"; E="  -0.25%  " },
    @{ Row=23; D="236.15"; E="  -0.14%  " },
    @{ Row=24; E="  +1.14%  " },
    @{ Row=25; D="2.45"; E="  +0.76%  " },
    @{ Row=26; E="  -0.05%  " },
    @{ Row=27; D="24.57"; E="  -0.73%  " },
    @{ Row=28; D="2.35"; E="  -0.38%  " },
    @{ Row=29; D="9.36"; E="  +1.86%  " },
    @{ Row=30; D="32.04"; E="  +1.47%  " },
    @{ Row=31; B="Kaspa"; C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D="0.115"; E="  +13.62%  " },
    @{ Row=32; B="FirstDigitalUSD"; C="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D="1.00"; E="  -0.03%  " },
    @{ Row=33; D="5.08"; E="  +1.34%  " },
    @{ Row=34; D="17.99"; E="  +3.46%  " },
    @{ Row=35; E="  +1.16%  " },
    @{ Row=36; D="130.38"; E="  +13.21%  " },
    @{ Row=37; D="1.83"; E="  -0.75%  " },
    @{ Row=38; E="  +3.17%  " },
    @{ Row=39; E="  -1.78%  " },
    @{ Row=40; E="  -2.84%  " },
    @{ Row=41; E="  -0.82%  " },
    @{ Row=42; D="21.24"; E="  -6.67%  " },
    @{ Row=43; D="1.931.76"; E="  -0.68%  " },
    @{ Row=44; E="  -1.16%  " },
    @{ Row=45; E="  +1.70%  " },
    @{ Row=46; D="2.78"; E="  +1.43%  " },
    @{ Row=47; D="9.26"; E="  -7.61%  " },
    @{ Row=48; D="2.601.01"; E="  +0.77%  " },
    @{ Row=49; E="  +2.38%  " },
    @{ Row=50; D="52.08"; E="  -2.28%  " },
    @{ Row=51; D="71.71"; E="  -0.80%  " }

)

foreach ($item in $updates) {
    $row = $item.Row

    if ($item.ContainsKey("B")) {
        $ws.Cells.Item($row, 2).Value = $item.B
    }
    if ($item.ContainsKey("C")) {
        $ws.Cells.Item($row, 3).Value = $item.C
    }
    if ($item.ContainsKey("D")) {
        # Price column holds text like "43.111.92" / "1.00" that Excel would
        # otherwise auto-convert to a number; format as Text first so the
        # literal string is preserved, exactly like the source data.
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $item.E
    }
}
